$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F17").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F19").Value = "ppe"
$ws.Range("F22").Value = "application instructions"
$ws.Range("F23").Value = "32_physical_and_chemical_hazards"
$ws.Range("F26").Value = "application instructions"
$ws.Range("F27").Value = "application instructions"
$ws.Range("F29").Value = "application instructions"
$ws.Range("F33").Value = "135_product_information"
$ws.Range("F34").Value = "use restrictions"
$ws.Range("F36").Value = "use restrictions || application instructions"
$ws.Range("F37").Value = "application instructions"
$ws.Range("F38").Value = "application instructions"
$ws.Range("F39").Value = "application instructions"
$ws.Range("F40").Value = "application instructions"
$ws.Range("F43").Value = "application instructions"
$ws.Range("F45").Value = "application instructions"
$ws.Range("F56").Value = "mixing"
$ws.Range("F57").Value = "mixing"
$ws.Range("F58").Value = "mixing"
$ws.Range("F59").Value = "mixing"
$ws.Range("F60").Value = "mixing"
$ws.Range("F61").Value = "mixing"
$ws.Range("F65").Value = "application instructions"
$ws.Range("F66").Value = "application instructions"
$ws.Range("F67").Value = "application instructions"
$ws.Range("F68").Value = "mixing"
$ws.Range("F69").Value = "mixing"
$ws.Range("F70").Value = "mixing"
$ws.Range("F72").Value = "mixing"
$ws.Range("F95").Value = "mixing"
$ws.Range("F96").Value = "mixing"
$ws.Range("F99").Value = "mixing"
$ws.Range("F101").Value = "mixing"
$ws.Range("F102").Value = "mixing"
$ws.Range("F108").Value = "mixing"
$ws.Range("F109").Value = "mixing"
$ws.Range("F110").Value = "mixing"
$ws.Range("F111").Value = "mixing"
$ws.Range("F112").Value = "mixing"
$ws.Range("F113").Value = "mixing"
$ws.Range("F116").Value = "mixing"
$ws.Range("F117").Value = "mixing"
$ws.Range("F118").Value = "mixing"
$ws.Range("F119").Value = "mixing"
$ws.Range("F121").Value = "mixing"
$ws.Range("F122").Value = "mixing"
$ws.Range("F123").Value = "mixing"
$ws.Range("F124").Value = "mixing"
$ws.Range("F125").Value = "mixing"
$ws.Range("F128").Value = "mixing"
$ws.Range("F134").Value = "mixing"
$ws.Range("F135").Value = "use restrictions"
$ws.Range("F136").Value = "mixing"
$ws.Range("F142").Value = "application instructions"
$ws.Range("F146").Value = "use restrictions"
$ws.Range("F147").Value = "use restrictions"
$ws.Range("F149").Value = "application instructions"
$ws.Range("F150").Value = "use restrictions"
$ws.Range("F156").Value = "use restrictions || application instructions"
$ws.Range("F157").Value = "use restrictions"
$ws.Range("F158").Value = "use restrictions"
$ws.Range("F159").Value = "use restrictions"
$ws.Range("F170").Value = "application instructions"
$ws.Range("F171").Value = "use restrictions"
$ws.Range("F173").Value = "application instructions"
$ws.Range("F174").Value = "application instructions"
$ws.Range("F175").Value = "application instructions"
$ws.Range("F176").Value = "application instructions"
$ws.Range("F177").Value = "application instructions"
$ws.Range("F178").Value = "application instructions"
$ws.Range("F179").Value = "application instructions"
$ws.Range("F182").Value = "safety procedures"
$ws.Range("F183").Value = "safety procedures"
$ws.Range("F184").Value = "safety procedures"
$ws.Range("F185").Value = "off target movement"
$ws.Range("F186").Value = "off target movement"
$ws.Range("F187").Value = "off target movement"
$ws.Range("F190").Value = "off target movement"
$ws.Range("F191").Value = "off target movement"
$ws.Range("F192").Value = "off target movement"
$ws.Range("F193").Value = "off target movement"
$ws.Range("F194").Value = "off target movement"
$ws.Range("F195").Value = "off target movement"
$ws.Range("F196").Value = "off target movement"
$ws.Range("F203").Value = "use restrictions"
$ws.Range("F204").Value = "use restrictions"
$ws.Range("F206").Value = "154_pesticide_storage"
